# Applies scheduled market-data refresh updates to the Pandaemonium_Profits sheets.
# For each affected Leve row: set new currentAveragePrice/LevePrice/LeveProfit values;
# clear cells that no longer carry a value, and set cells that are newly populated.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 773.5714
$ws.Range("I28").Value = 220.90909
$ws.Range("J28").Value = 2800
$ws.Range("K28").Value = 220.90909
$ws.Range("L28").Value = 2800
$ws.Range("M28").Value = 264.09091
$ws.Range("N28").Value = -3770

$ws.Range("H51").Value = 1500
$ws.Range("J51").Value = 1500
$ws.Range("L51").Value = 1500
$ws.Range("N51").Value = -2468

$ws.Range("H69").Value = 5944.421
$ws.Range("I69").Value = 6411.7
$ws.Range("J69").Value = 5425.222
$ws.Range("K69").Value = 19235.1
$ws.Range("L69").Value = 16275.666
$ws.Range("M69").Value = -18361.1
$ws.Range("N69").Value = -18023.666

$ws.Range("H72").Value = 5944.421
$ws.Range("I72").Value = 6411.7
$ws.Range("J72").Value = 5425.222
$ws.Range("K72").Value = 57705.3
$ws.Range("L72").Value = 48826.998
$ws.Range("M72").Value = -53337.3
$ws.Range("N72").Value = -57562.998

$ws.Range("H98").Value = 654.8125
$ws.Range("I98").Value = 654.8125
$ws.Range("K98").Value = 654.8125
$ws.Range("M98").Value = 843.1875

$ws.Range("H107").Value = 857.55554
$ws.Range("I107").Value = 717.0833
$ws.Range("J107").Value = 1138.5
$ws.Range("K107").Value = 717.0833
$ws.Range("L107").Value = 1138.5
$ws.Range("M107").Value = 1202.9167
$ws.Range("N107").Value = -4978.5

$ws.Range("H113").Value = 3781.111
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 4406
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 4406
$ws.Range("M113").Value = 254
$ws.Range("N113").Value = -10914

$ws.Range("H118").Value = 696.6667
$ws.Range("I118").Value = 595
$ws.Range("K118").Value = 1785
$ws.Range("M118").Value = -128

$ws.Range("H122").Value = 654.8125
$ws.Range("I122").Value = 654.8125
$ws.Range("K122").Value = 1964.4375
$ws.Range("M122").Value = 485.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7247.8936
$ws.Range("I32").Value = 5955.976
$ws.Range("J32").Value = 18100
$ws.Range("K32").Value = 5955.976
$ws.Range("L32").Value = 18100
$ws.Range("M32").Value = -5668.976
$ws.Range("N32").Value = -18674

$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H110").Value = 833.3333
$ws.Range("I110").Value = 833.3333
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 833.3333
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1211.6667
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3484
$ws.Range("I86").Value = 3484
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3484
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2361
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 3484
$ws.Range("I89").Value = 3484
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 17420
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -11804
$ws.Range("N89").Value = -11804

$ws.Range("H132").Value = 63498
$ws.Range("J132").Value = 63498
$ws.Range("L132").Value = 63498
$ws.Range("N132").Value = -73618

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1501.25
$ws.Range("I16").Value = 1202
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1202
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -915
$ws.Range("N16").Value = -2574

$ws.Range("H31").Value = 3832.3635
$ws.Range("I31").Value = 2884.6924
$ws.Range("K31").Value = 2884.6924
$ws.Range("M31").Value = -2589.6924

$ws.Range("H34").Value = 3832.3635
$ws.Range("I34").Value = 2884.6924
$ws.Range("K34").Value = 2884.6924
$ws.Range("M34").Value = -2682.6924

$ws.Range("H107").Value = 916.9
$ws.Range("I107").Value = 904.75
$ws.Range("J107").Value = 941.2
$ws.Range("K107").Value = 904.75
$ws.Range("L107").Value = 941.2
$ws.Range("M107").Value = 1015.25
$ws.Range("N107").Value = -4781.2

$ws.Range("H113").Value = 1501.25
$ws.Range("I113").Value = 1202
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1202
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 968
$ws.Range("N113").Value = -6340

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5491.4
$ws.Range("I70").Value = 5379.8667
$ws.Range("J70").Value = 5521.8184
$ws.Range("K70").Value = 5379.8667
$ws.Range("L70").Value = 5521.8184
$ws.Range("M70").Value = -5109.8667
$ws.Range("N70").Value = -6061.8184

$ws.Range("H73").Value = 5491.4
$ws.Range("I73").Value = 5379.8667
$ws.Range("J73").Value = 5521.8184
$ws.Range("K73").Value = 5379.8667
$ws.Range("L73").Value = 5521.8184
$ws.Range("M73").Value = -4443.8667
$ws.Range("N73").Value = -7393.8184

$ws.Range("H107").Value = 627.6
$ws.Range("I107").Value = 438.57144
$ws.Range("J107").Value = 1068.6666
$ws.Range("K107").Value = 438.57144
$ws.Range("L107").Value = 1068.6666
$ws.Range("M107").Value = 1481.42856
$ws.Range("N107").Value = -4908.6666

$ws.Range("H113").Value = 2364
$ws.Range("I113").Value = 3374.75
$ws.Range("J113").Value = 1353.25
$ws.Range("K113").Value = 3374.75
$ws.Range("L113").Value = 1353.25
$ws.Range("M113").Value = -1204.75
$ws.Range("N113").Value = -5693.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 595.0625
$ws.Range("I22").Value = 574.2
$ws.Range("J22").Value = 604.5454999999999
$ws.Range("K22").Value = 574.2
$ws.Range("L22").Value = 604.5454999999999
$ws.Range("M22").Value = -279.2
$ws.Range("N22").Value = -1194.5455

$ws.Range("H27").Value = 595.0625
$ws.Range("I27").Value = 574.2
$ws.Range("J27").Value = 604.5454999999999
$ws.Range("K27").Value = 574.2
$ws.Range("L27").Value = 604.5454999999999
$ws.Range("M27").Value = -467.2
$ws.Range("N27").Value = -818.5454999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1337.2609
$ws.Range("I107").Value = 603.4666999999999
$ws.Range("J107").Value = 2713.125
$ws.Range("K107").Value = 1810.4001
$ws.Range("L107").Value = 8139.375
$ws.Range("M107").Value = 109.5999000000002
$ws.Range("N107").Value = -11979.375

$ws.Range("H113").Value = 982.8421
$ws.Range("I113").Value = 346.42856
$ws.Range("K113").Value = 1039.28568
$ws.Range("M113").Value = 1130.71432
